# The commit swaps the two embedded DrawingML themes used by this deck:
#   ppt/theme/theme1.xml (bound to the slide master, i.e. the deck's
#   visible design) switches from the custom "Integral" / "Red Violet"
#   colour scheme to the stock PowerPoint "Office Theme" / "Office"
#   colour scheme, while ppt/theme/theme2.xml (bound to the notes
#   master) takes on the old "Integral" / "Red Violet" colours.
#
# PowerPoint's automation model doesn't expose a way to rebind which
# OOXML theme part backs a given master, but it does expose each
# master's resolved 12-slot ThemeColorScheme as settable RGB values, so
# we push the standard "Office Theme" palette onto every theme surface
# PowerPoint's object model lets us reach (the presentation's slide
# master design and its notes master) to mirror the content change as
# closely as the object model allows.

$p = $ppt.ActivePresentation

# Standard "Office Theme" colour scheme (from the target theme XML),
# expressed as OLE RGB() values (0xBBGGRR) the way ThemeColor.RGB wants
# them.
$officeColors = @{
    1  = 0         # dk1      000000
    2  = 16777215  # lt1      FFFFFF
    3  = 6968388   # dk2      44546A
    4  = 15132391  # lt2      E7E6E6
    5  = 13998939  # accent1  5B9BD5
    6  = 3243501   # accent2  ED7D31
    7  = 10855845  # accent3  A5A5A5
    8  = 49407     # accent4  FFC000
    9  = 12874308  # accent5  4472C4
    10 = 4697456   # accent6  70AD47
    11 = 12673797  # hlink    0563C1
    12 = 7491477   # folHlink 954F72
}

function Set-OfficeColors($colorScheme) {
    if ($colorScheme -eq $null) { return }
    for ($i = 1; $i -le 12; $i++) {
        $colorScheme.Item($i).RGB = $officeColors[$i]
    }
}

# Slide (design) master theme -> ppt/theme/theme1.xml
try {
    $master = $p.SlideMaster
    Set-OfficeColors $master.Theme.ThemeColorScheme
} catch {
    Write-Host "SlideMaster theme update failed:" $_.Exception.Message
}

# Notes master theme -> ppt/theme/theme2.xml (where the object model
# resolves it distinctly from the slide master)
try {
    $notesMaster = $p.NotesMaster
    Set-OfficeColors $notesMaster.Theme.ThemeColorScheme
} catch {
    Write-Host "NotesMaster theme update failed:" $_.Exception.Message
}
